$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.281134"
$ws.Range("H2").Value = [double]"3.843402"
$ws.Range("I2").Value = [double]"0.007312702338676299"
$ws.Range("J2").Value = [double]"0.007312702338676299"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"3.151158666666667"
$ws.Range("N2").Value = [double]"9.453476"
$ws.Range("O2").Value = [double]"0.03114707555614071"
$ws.Range("P2").Value = [double]"0.03114707555614071"
$ws.Range("Q2").Value = [double]"4.037056507261334"
$ws.Range("R2").Value = [double]"36.33350856535201"
$ws.Range("S2").Value = [double]"0.0002277692922623176"
$ws.Range("T2").Value = [double]"0.0002277692922623175"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.281134"
$ws.Range("H3").Value = [double]"3.843402"
$ws.Range("I3").Value = [double]"0.007312702338676299"
$ws.Range("J3").Value = [double]"0.007312702338676299"
$ws.Range("M3").Value = [double]"5.038243666666667"
$ws.Range("O3").Value = [double]"0.04979963650066307"
$ws.Range("P3").Value = [double]"0.04979963650066306"
$ws.Range("Q3").Value = [double]"6.454665261651334"
$ws.Range("R3").Value = [double]"58.09198735486201"
$ws.Range("S3").Value = [double]"0.0003641699183036284"
$ws.Range("T3").Value = [double]"0.0003641699183036284"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.281134"
$ws.Range("H4").Value = [double]"3.843402"
$ws.Range("I4").Value = [double]"0.007312702338676299"
$ws.Range("J4").Value = [double]"0.007312702338676299"
$ws.Range("M4").Value = [double]"92.91163899999999"
$ws.Range("N4").Value = [double]"278.734917"
$ws.Range("O4").Value = [double]"0.9183688116343246"
$ws.Range("P4").Value = [double]"0.9183688116343246"
$ws.Range("Q4").Value = [double]"119.032259718626"
$ws.Range("R4").Value = [double]"1071.290337467634"
$ws.Range("S4").Value = [double]"0.006715757756605698"
$ws.Range("T4").Value = [double]"0.006715757756605698"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"1.281134"
$ws.Range("H5").Value = [double]"3.843402"
$ws.Range("I5").Value = [double]"0.007312702338676299"
$ws.Range("J5").Value = [double]"0.007312702338676299"
$ws.Range("M5").Value = [double]"0.06924866666666667"
$ws.Range("N5").Value = [double]"0.207746"
$ws.Range("O5").Value = [double]"0.0006844763088715736"
$ws.Range("P5").Value = [double]"0.0006844763088715734"
$ws.Range("Q5").Value = [double]"0.08871682132133334"
$ws.Range("R5").Value = [double]"0.798451391892"
$ws.Range("S5").Value = [double]"5.005371504653677E-06"
$ws.Range("T5").Value = [double]"5.005371504653676E-06"
$ws.Range("I6").Value = [double]"0.9398544320918915"
$ws.Range("J6").Value = [double]"0.9398544320918915"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"3.151158666666667"
$ws.Range("N6").Value = [double]"9.453476"
$ws.Range("O6").Value = [double]"0.03114707555614071"
$ws.Range("P6").Value = [double]"0.03114707555614071"
$ws.Range("Q6").Value = [double]"518.8568158842613"
$ws.Range("R6").Value = [double]"4669.711342958352"
$ws.Range("S6").Value = [double]"0.02927371700813986"
$ws.Range("T6").Value = [double]"0.02927371700813986"
$ws.Range("I7").Value = [double]"0.9398544320918915"
$ws.Range("J7").Value = [double]"0.9398544320918915"
$ws.Range("M7").Value = [double]"5.038243666666667"
$ws.Range("O7").Value = [double]"0.04979963650066307"
$ws.Range("P7").Value = [double]"0.04979963650066306"
$ws.Range("Q7").Value = [double]"829.5764647424014"
$ws.Range("R7").Value = [double]"7466.188182681613"
$ws.Range("S7").Value = [double]"0.04680440908171331"
$ws.Range("T7").Value = [double]"0.04680440908171331"
$ws.Range("I8").Value = [double]"0.9398544320918915"
$ws.Range("J8").Value = [double]"0.9398544320918915"
$ws.Range("M8").Value = [double]"92.91163899999999"
$ws.Range("N8").Value = [double]"278.734917"
$ws.Range("O8").Value = [double]"0.9183688116343246"
$ws.Range("P8").Value = [double]"0.9183688116343246"
$ws.Range("Q8").Value = [double]"15298.44805343387"
$ws.Range("R8").Value = [double]"137686.0324809049"
$ws.Range("S8").Value = [double]"0.8631329979094834"
$ws.Range("T8").Value = [double]"0.8631329979094834"
$ws.Range("I9").Value = [double]"0.9398544320918915"
$ws.Range("J9").Value = [double]"0.9398544320918915"
$ws.Range("M9").Value = [double]"0.06924866666666667"
$ws.Range("N9").Value = [double]"0.207746"
$ws.Range("O9").Value = [double]"0.0006844763088715736"
$ws.Range("P9").Value = [double]"0.0006844763088715734"
$ws.Range("Q9").Value = [double]"11.40220042582133"
$ws.Range("R9").Value = [double]"102.619803832392"
$ws.Range("S9").Value = [double]"0.0006433080925548468"
$ws.Range("T9").Value = [double]"0.0006433080925548467"
$ws.Range("G10").Value = [double]"9.213772333333333"
$ws.Range("H10").Value = [double]"27.641317"
$ws.Range("I10").Value = [double]"0.05259213672418158"
$ws.Range("J10").Value = [double]"0.05259213672418158"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"3.151158666666667"
$ws.Range("N10").Value = [double]"9.453476"
$ws.Range("O10").Value = [double]"0.03114707555614071"
$ws.Range("P10").Value = [double]"0.03114707555614071"
$ws.Range("Q10").Value = [double]"29.03405854087689"
$ws.Range("R10").Value = [double]"261.306526867892"
$ws.Range("S10").Value = [double]"0.001638091256206966"
$ws.Range("T10").Value = [double]"0.001638091256206966"
$ws.Range("G11").Value = [double]"9.213772333333333"
$ws.Range("H11").Value = [double]"27.641317"
$ws.Range("I11").Value = [double]"0.05259213672418158"
$ws.Range("J11").Value = [double]"0.05259213672418158"
$ws.Range("M11").Value = [double]"5.038243666666667"
$ws.Range("O11").Value = [double]"0.04979963650066307"
$ws.Range("P11").Value = [double]"0.04979963650066306"
$ws.Range("Q11").Value = [double]"46.42123010452523"
$ws.Range("R11").Value = [double]"417.7910709407271"
$ws.Range("S11").Value = [double]"0.002619069291657416"
$ws.Range("T11").Value = [double]"0.002619069291657415"
$ws.Range("G12").Value = [double]"9.213772333333333"
$ws.Range("H12").Value = [double]"27.641317"
$ws.Range("I12").Value = [double]"0.05259213672418158"
$ws.Range("J12").Value = [double]"0.05259213672418158"
$ws.Range("M12").Value = [double]"92.91163899999999"
$ws.Range("N12").Value = [double]"278.734917"
$ws.Range("O12").Value = [double]"0.9183688116343246"
$ws.Range("P12").Value = [double]"0.9183688116343246"
$ws.Range("Q12").Value = [double]"856.0666888628542"
$ws.Range("R12").Value = [double]"7704.600199765689"
$ws.Range("S12").Value = [double]"0.04829897810469656"
$ws.Range("T12").Value = [double]"0.04829897810469656"
$ws.Range("G13").Value = [double]"9.213772333333333"
$ws.Range("H13").Value = [double]"27.641317"
$ws.Range("I13").Value = [double]"0.05259213672418158"
$ws.Range("J13").Value = [double]"0.05259213672418158"
$ws.Range("M13").Value = [double]"0.06924866666666667"
$ws.Range("N13").Value = [double]"0.207746"
$ws.Range("O13").Value = [double]"0.0006844763088715736"
$ws.Range("P13").Value = [double]"0.0006844763088715734"
$ws.Range("Q13").Value = [double]"0.6380414490535555"
$ws.Range("R13").Value = [double]"5.742373041482"
$ws.Range("S13").Value = [double]"3.599807162063693E-05"
$ws.Range("T13").Value = [double]"3.599807162063693E-05"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.042174"
$ws.Range("H14").Value = [double]"0.126522"
$ws.Range("I14").Value = [double]"0.0002407288452506406"
$ws.Range("J14").Value = [double]"0.0002407288452506406"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"3.151158666666667"
$ws.Range("N14").Value = [double]"9.453476"
$ws.Range("O14").Value = [double]"0.03114707555614071"
$ws.Range("P14").Value = [double]"0.03114707555614071"
$ws.Range("Q14").Value = [double]"0.132896965608"
$ws.Range("R14").Value = [double]"1.196072690472"
$ws.Range("S14").Value = [double]"7.497999531564209E-06"
$ws.Range("T14").Value = [double]"7.497999531564207E-06"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.042174"
$ws.Range("H15").Value = [double]"0.126522"
$ws.Range("I15").Value = [double]"0.0002407288452506406"
$ws.Range("J15").Value = [double]"0.0002407288452506406"
$ws.Range("M15").Value = [double]"5.038243666666667"
$ws.Range("O15").Value = [double]"0.04979963650066307"
$ws.Range("P15").Value = [double]"0.04979963650066306"
$ws.Range("Q15").Value = [double]"0.212482888398"
$ws.Range("R15").Value = [double]"1.912345995582"
$ws.Range("S15").Value = [double]"1.198820898870627E-05"
$ws.Range("T15").Value = [double]"1.198820898870627E-05"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.042174"
$ws.Range("H16").Value = [double]"0.126522"
$ws.Range("I16").Value = [double]"0.0002407288452506406"
$ws.Range("J16").Value = [double]"0.0002407288452506406"
$ws.Range("M16").Value = [double]"92.91163899999999"
$ws.Range("N16").Value = [double]"278.734917"
$ws.Range("O16").Value = [double]"0.9183688116343246"
$ws.Range("P16").Value = [double]"0.9183688116343246"
$ws.Range("Q16").Value = [double]"3.918455463186"
$ws.Range("R16").Value = [double]"35.266099168674"
$ws.Range("S16").Value = [double]"0.000221077863538934"
$ws.Range("T16").Value = [double]"0.000221077863538934"
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.042174"
$ws.Range("H17").Value = [double]"0.126522"
$ws.Range("I17").Value = [double]"0.0002407288452506406"
$ws.Range("J17").Value = [double]"0.0002407288452506406"
$ws.Range("M17").Value = [double]"0.06924866666666667"
$ws.Range("N17").Value = [double]"0.207746"
$ws.Range("O17").Value = [double]"0.0006844763088715736"
$ws.Range("P17").Value = [double]"0.0006844763088715734"
$ws.Range("Q17").Value = [double]"0.002920493268"
$ws.Range("R17").Value = [double]"0.026284439412"
$ws.Range("S17").Value = [double]"1.647731914360747E-07"
$ws.Range("T17").Value = [double]"1.647731914360747E-07"
